# Updated test data for FC test cases
# - Rename the "Battery Alarm (A)" / "Battery Standby (A)" labels used on the
#   "Add Panels" sheet to "Alarm Current(A)" / "Standby Current(A)".
# - Leave the active selection on cell O8 (last selected cell before save).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# These labels appear twice on the sheet (row 2 header pair and row 8 data pair)
# and both point at the same shared-string entries, so update every occurrence.
$ws.Range("F2").Value = "Alarm Current(A)"
$ws.Range("G2").Value = "Standby Current(A)"
$ws.Range("N8").Value = "Alarm Current(A)"
$ws.Range("O8").Value = "Standby Current(A)"

# Restore the saved selection state (active cell O8 on the "Add Panels" sheet).
$ws.Activate()
$ws.Range("O8").Select()
